$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 456.8
$ws.Range("I33").Value = 445.64285
$ws.Range("K33").Value = 445.64285
$ws.Range("M33").Value = -216.64285
$ws.Range("H64").Value = 355406.34
$ws.Range("I64").Value = 603528.3
$ws.Range("J64").Value = 3900.25
$ws.Range("K64").Value = 603528.3
$ws.Range("L64").Value = 3900.25
$ws.Range("M64").Value = -603280.3
$ws.Range("N64").Value = -4396.25
$ws.Range("H67").Value = 355406.34
$ws.Range("I67").Value = 603528.3
$ws.Range("J67").Value = 3900.25
$ws.Range("K67").Value = 603528.3
$ws.Range("L67").Value = 3900.25
$ws.Range("M67").Value = -602670.3
$ws.Range("N67").Value = -5616.25
$ws.Range("H100").Value = 2768.5715
$ws.Range("I100").Value = 1933.3334
$ws.Range("J100").Value = 2996.3635
$ws.Range("K100").Value = 1933.3334
$ws.Range("L100").Value = 2996.3635
$ws.Range("M100").Value = -1392.3334
$ws.Range("N100").Value = -4078.3635
$ws.Range("H135").Value = 48388988
$ws.Range("I135").Value = 19232330
$ws.Range("J135").Value = 200003600
$ws.Range("K135").Value = 173090970
$ws.Range("L135").Value = 1800032400
$ws.Range("M135").Value = -173088435
$ws.Range("N135").Value = -1800037470
$ws.Range("H138").Value = 3766.4285
$ws.Range("J138").Value = 4230.965
$ws.Range("L138").Value = 12692.895
$ws.Range("N138").Value = -22972.895

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1757.138
$ws.Range("I2").Value = 1864
$ws.Range("J2").Value = 1476.625
$ws.Range("K2").Value = 1864
$ws.Range("L2").Value = 1476.625
$ws.Range("M2").Value = -1751
$ws.Range("N2").Value = -1702.625
$ws.Range("H110").Value = 1516.2307
$ws.Range("I110").Value = 1446.4546
$ws.Range("J110").Value = 1900
$ws.Range("K110").Value = 1446.4546
$ws.Range("L110").Value = 1900
$ws.Range("M110").Value = 598.5454
$ws.Range("N110").Value = -5990
$ws.Range("H116").Value = 1757.138
$ws.Range("I116").Value = 1864
$ws.Range("J116").Value = 1476.625
$ws.Range("K116").Value = 1864
$ws.Range("L116").Value = 1476.625
$ws.Range("M116").Value = 430
$ws.Range("N116").Value = -6064.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1757.138
$ws.Range("I3").Value = 1864
$ws.Range("J3").Value = 1476.625
$ws.Range("K3").Value = 1864
$ws.Range("L3").Value = 1476.625
$ws.Range("M3").Value = -1750
$ws.Range("N3").Value = -1704.625
$ws.Range("H24").Value = 3086.1428
$ws.Range("I24").Value = 1724.6
$ws.Range("J24").Value = 6490
$ws.Range("K24").Value = 1724.6
$ws.Range("L24").Value = 6490
$ws.Range("M24").Value = -1489.6
$ws.Range("N24").Value = -6960
$ws.Range("H99").Value = 2127.7778
$ws.Range("I99").Value = 1953.3334
$ws.Range("K99").Value = 1953.3334
$ws.Range("M99").Value = -455.3334
$ws.Range("H105").Value = 5592.36
$ws.Range("I105").Value = 5263.5454
$ws.Range("J105").Value = 8003.6665
$ws.Range("K105").Value = 5263.5454
$ws.Range("L105").Value = 8003.6665
$ws.Range("M105").Value = -3516.5454
$ws.Range("N105").Value = -11497.6665
$ws.Range("H107").Value = 2530
$ws.Range("H111").Value = 59650
$ws.Range("J111").Value = 59650
$ws.Range("L111").Value = 59650
$ws.Range("N111").Value = -67830
$ws.Range("H134").Value = 27483.875
$ws.Range("I134").Value = 2466.0857
$ws.Range("J134").Value = 202608.4
$ws.Range("K134").Value = 7398.257100000001
$ws.Range("L134").Value = 607825.2
$ws.Range("M134").Value = -4863.257100000001
$ws.Range("N134").Value = -612895.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3355.9092
$ws.Range("I62").Value = 3341.5
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 3341.5
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -2717.5
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 3355.9092
$ws.Range("I65").Value = 3341.5
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 16707.5
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -13587.5
$ws.Range("N65").Value = -23740
$ws.Range("H102").Value = 45000
$ws.Range("J102").Value = 45000
$ws.Range("L102").Value = 45000
$ws.Range("N102").Value = -49868
$ws.Range("H107").Value = 1794.6666
$ws.Range("I107").Value = 2215.7144
$ws.Range("J107").Value = 1205.2
$ws.Range("K107").Value = 2215.7144
$ws.Range("L107").Value = 1205.2
$ws.Range("M107").Value = -295.7143999999998
$ws.Range("N107").Value = -5045.2
$ws.Range("H122").Value = 10199.32
$ws.Range("I122").Value = 6525
$ws.Range("J122").Value = 11928.412
$ws.Range("K122").Value = 19575
$ws.Range("L122").Value = 35785.236
$ws.Range("M122").Value = -17125
$ws.Range("N122").Value = -40685.236
$ws.Range("H134").Value = 2761.1343
$ws.Range("I134").Value = 1547.6487
$ws.Range("J134").Value = 4257.7666
$ws.Range("K134").Value = 4642.9461
$ws.Range("L134").Value = 12773.2998
$ws.Range("M134").Value = -2107.9461
$ws.Range("N134").Value = -17843.2998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5151.6523
$ws.Range("I68").Value = 648.75
$ws.Range("J68").Value = 7553.2
$ws.Range("K68").Value = 1946.25
$ws.Range("L68").Value = 22659.6
$ws.Range("M68").Value = -1135.25
$ws.Range("N68").Value = -24281.6
$ws.Range("H71").Value = 5151.6523
$ws.Range("I71").Value = 648.75
$ws.Range("J71").Value = 7553.2
$ws.Range("K71").Value = 5838.75
$ws.Range("L71").Value = 67978.8
$ws.Range("M71").Value = -1782.75
$ws.Range("N71").Value = -76090.8
$ws.Range("H108").Value = 3343
$ws.Range("I108").Value = 1499.5
$ws.Range("K108").Value = 4498.5
$ws.Range("M108").Value = -1618.5
$ws.Range("H132").Value = 1958.375
$ws.Range("I132").Value = 2106.75
$ws.Range("J132").Value = 1810
$ws.Range("K132").Value = 18960.75
$ws.Range("L132").Value = 16290
$ws.Range("M132").Value = -16430.75
$ws.Range("N132").Value = -21350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 10416.909
$ws.Range("J18").Value = 10509.556
$ws.Range("L18").Value = 10509.556
$ws.Range("N18").Value = -11095.556
$ws.Range("H46").Value = 16908.6
$ws.Range("I46").Value = 10020
$ws.Range("J46").Value = 18630.75
$ws.Range("K46").Value = 10020
$ws.Range("L46").Value = 18630.75
$ws.Range("M46").Value = -9864
$ws.Range("N46").Value = -18942.75
$ws.Range("H57").Value = 15500
$ws.Range("I57").Value = 3000
$ws.Range("J57").Value = 19666.666
$ws.Range("K57").Value = 3000
$ws.Range("L57").Value = 19666.666
$ws.Range("M57").Value = -2180
$ws.Range("N57").Value = -21306.666
$ws.Range("H70").Value = 5289.962
$ws.Range("I70").Value = 4834.923
$ws.Range("J70").Value = 5513.1885
$ws.Range("K70").Value = 4834.923
$ws.Range("L70").Value = 5513.1885
$ws.Range("M70").Value = -4564.923
$ws.Range("N70").Value = -6053.1885
$ws.Range("H73").Value = 5289.962
$ws.Range("I73").Value = 4834.923
$ws.Range("J73").Value = 5513.1885
$ws.Range("K73").Value = 4834.923
$ws.Range("L73").Value = 5513.1885
$ws.Range("M73").Value = -3898.923
$ws.Range("N73").Value = -7385.1885
$ws.Range("H80").Value = 7802.909
$ws.Range("I80").Value = 18616.334
$ws.Range("J80").Value = 3747.875
$ws.Range("K80").Value = 18616.334
$ws.Range("L80").Value = 3747.875
$ws.Range("M80").Value = -17618.334
$ws.Range("N80").Value = -5743.875
$ws.Range("H83").Value = 7802.909
$ws.Range("I83").Value = 18616.334
$ws.Range("J83").Value = 3747.875
$ws.Range("K83").Value = 93081.67
$ws.Range("L83").Value = 18739.375
$ws.Range("M83").Value = -88089.67
$ws.Range("N83").Value = -28723.375
$ws.Range("H102").Value = 3259.7021
$ws.Range("I102").Value = 2939.2646
$ws.Range("J102").Value = 4097.769
$ws.Range("K102").Value = 2939.2646
$ws.Range("L102").Value = 4097.769
$ws.Range("M102").Value = -1317.2646
$ws.Range("N102").Value = -7341.769
$ws.Range("H122").Value = 6838
$ws.Range("I122").Value = 11120
$ws.Range("J122").Value = 3269.6667
$ws.Range("K122").Value = 33360
$ws.Range("L122").Value = 9809.000100000001
$ws.Range("M122").Value = -30910
$ws.Range("N122").Value = -14709.0001
$ws.Range("H126").Value = 3146.95
$ws.Range("I126").Value = 2028.5714
$ws.Range("J126").Value = 3749.1538
$ws.Range("K126").Value = 6085.7142
$ws.Range("L126").Value = 11247.4614
$ws.Range("M126").Value = -3615.7142
$ws.Range("N126").Value = -16187.4614
$ws.Range("H132").Value = 32217.082
$ws.Range("I132").Value = 57178.445
$ws.Range("J132").Value = 8569.474
$ws.Range("K132").Value = 171535.335
$ws.Range("L132").Value = 25708.422
$ws.Range("M132").Value = -169005.335
$ws.Range("N132").Value = -30768.422

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 4966.6665
$ws.Range("J59").Value = 4966.6665
$ws.Range("L59").Value = 4966.6665
$ws.Range("N59").Value = -6274.6665
$ws.Range("H111").Value = 45193.5
$ws.Range("J111").Value = 45193.5
$ws.Range("L111").Value = 45193.5
$ws.Range("N111").Value = -53373.5
$ws.Range("H122").Value = 8357.429
$ws.Range("I122").Value = 5600.5713
$ws.Range("J122").Value = 11114.286
$ws.Range("K122").Value = 16801.7139
$ws.Range("L122").Value = 33342.858
$ws.Range("M122").Value = -14351.7139
$ws.Range("N122").Value = -38242.858
$ws.Range("H132").Value = 3073.861
$ws.Range("I132").Value = 2774.7778
$ws.Range("J132").Value = 3971.111
$ws.Range("K132").Value = 8324.3334
$ws.Range("L132").Value = 11913.333
$ws.Range("M132").Value = -5794.3334
$ws.Range("N132").Value = -16973.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").ClearContents()
$ws.Range("N114").Value = 0
$ws.Range("H122").Value = 4343.407
$ws.Range("I122").Value = 2404.3845
$ws.Range("J122").Value = 6143.9287
$ws.Range("K122").Value = 7213.1535
$ws.Range("L122").Value = 18431.7861
$ws.Range("M122").Value = -4763.1535
$ws.Range("N122").Value = -23331.7861
$ws.Range("H126").Value = 1347.1177
$ws.Range("I126").Value = 1354.3846
$ws.Range("J126").Value = 1323.5
$ws.Range("K126").Value = 4063.1538
$ws.Range("L126").Value = 3970.5
$ws.Range("M126").Value = -1593.1538
$ws.Range("N126").Value = -8910.5
$ws.Range("H132").Value = 2857.1462
$ws.Range("I132").Value = 2117.2424
$ws.Range("J132").Value = 5909.25
$ws.Range("K132").Value = 6351.7272
$ws.Range("L132").Value = 17727.75
$ws.Range("M132").Value = -3821.7272
$ws.Range("N132").Value = -22787.75
